# Update "想去人数" (interest count) and a couple of "最低票价" (min price)
# values across the "展览", "演出" and "全部类型" worksheets to match the
# newly generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 152
$ws1.Range("F3").Value = 199
$ws1.Range("F6").Value = 1320
$ws1.Range("F7").Value = 71
$ws1.Range("F9").Value = 393
$ws1.Range("F11").Value = 816
$ws1.Range("F12").Value = 213
$ws1.Range("F14").Value = 317
$ws1.Range("F15").Value = 470
$ws1.Range("F16").Value = 92
$ws1.Range("F17").Value = 1052
$ws1.Range("F19").Value = 289
$ws1.Range("F20").Value = 412
$ws1.Range("F21").Value = 104
$ws1.Range("F22").Value = 223
$ws1.Range("F25").Value = 493
$ws1.Range("F26").Value = 445
$ws1.Range("F27").Value = 296

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = 150
$ws2.Range("F4").Value = 380
$ws2.Range("F6").Value = 48
$ws2.Range("F7").Value = 294
$ws2.Range("F12").Value = 146

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 150
$ws4.Range("F4").Value = 152
$ws4.Range("F5").Value = 199
$ws4.Range("F8").Value = 1320
$ws4.Range("F10").Value = 71
$ws4.Range("F11").Value = 380
$ws4.Range("F14").Value = 393
$ws4.Range("F15").Value = 48
$ws4.Range("F16").Value = 294
$ws4.Range("F18").Value = 816
$ws4.Range("F19").Value = 213
$ws4.Range("F21").Value = 317
$ws4.Range("F22").Value = 470
$ws4.Range("F23").Value = 92
$ws4.Range("F24").Value = 1052
$ws4.Range("F28").Value = 289
$ws4.Range("F29").Value = 412
$ws4.Range("F31").Value = 104
$ws4.Range("F33").Value = 223
$ws4.Range("F36").Value = 146
$ws4.Range("F38").Value = 493
$ws4.Range("F41").Value = 445
$ws4.Range("F42").Value = 296
